# Review more Coursera skills
# Adds 28 new Skill/Parent rows (rows 60-87) to the "Coursera" sheet,
# pushing the previously-last data row (row 99) further down the sheet
# (row 99 itself is untouched/unshifted - the sheet already had a gap).
#
# Cell writes are ordered to faithfully reproduce the original authoring
# sequence (column A filled top-to-bottom per row together with column B,
# except "Pricing Policy"'s Parent ("Policy") which was filled in later,
# right after the "Google Chart API" / "Google Services" row was entered).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coursera")

$ws.Range("A60").Value = "News Writing"
$ws.Range("B60").Value = "News"
$ws.Range("A61").Value = "Financial Data Analysis"
$ws.Range("B61").Value = "Finance"
$ws.Range("A62").Value = "Environmental Data Analysis"
$ws.Range("B62").Value = "Environmental Studies"
$ws.Range("A63").Value = "Pricing Policy"
$ws.Range("A64").Value = "Human Resource"
$ws.Range("B64").Value = "Business"
$ws.Range("A65").Value = "Performance Management"
$ws.Range("B65").Value = "Management"
$ws.Range("A66").Value = "Sociological Concept"
$ws.Range("B66").Value = "Sociology"
$ws.Range("A67").Value = "Economic Analysis"
$ws.Range("B67").Value = "Economics"
$ws.Range("A68").Value = "Economic Models"
$ws.Range("B68").Value = "Economics"
$ws.Range("A69").Value = "Product Management"
$ws.Range("B69").Value = "Management"
$ws.Range("A70").Value = "Statistical Analysis"
$ws.Range("B70").Value = "Statistics"
$ws.Range("A71").Value = "Business Analysis"
$ws.Range("B71").Value = "Analysis"
$ws.Range("A72").Value = "Regression Analysis"
$ws.Range("B72").Value = "Analysis"
$ws.Range("A73").Value = "Intellectual Property Law"
$ws.Range("B73").Value = "Intellectual Property"
$ws.Range("A74").Value = "Google Chart API"
$ws.Range("B74").Value = "Google Services"
$ws.Range("B63").Value = "Policy"
$ws.Range("A75").Value = "Uniform Resource Identifier (URI) Scheme"
$ws.Range("B75").Value = "World Wide Web (WWW)"
$ws.Range("A76").Value = "Financial Risk Modeling"
$ws.Range("B76").Value = "Finance"
$ws.Range("A77").Value = "RStudio"
$ws.Range("B77").Value = "R"
$ws.Range("A78").Value = "Auditor's Report"
$ws.Range("B78").Value = "Audit"
$ws.Range("A79").Value = "Data Reporting"
$ws.Range("B79").Value = "Data"
$ws.Range("A80").Value = "Sales Presentation"
$ws.Range("B80").Value = "Presentation"
$ws.Range("A81").Value = "Sales Management"
$ws.Range("B81").Value = "Sales"
$ws.Range("A82").Value = "Software Stack"
$ws.Range("B82").Value = "Information Technology"
$ws.Range("A83").Value = "Market Research"
$ws.Range("B83").Value = "Market"
$ws.Range("A84").Value = "Business Technology Management"
$ws.Range("B84").Value = "Business"
$ws.Range("A85").Value = "Marketing Research"
$ws.Range("B85").Value = "Marketing"
$ws.Range("A86").Value = "Market Analysis"
$ws.Range("B86").Value = "Market"
$ws.Range("A87").Value = "Cloud Computing Security"
$ws.Range("B87").Value = "Cloud Computing"

# Leave the cursor/selection where data entry finished (next empty row).
$ws.Range("A88").Select()
